# Regenerate save_data column G (K) with recalculated strike/touch counts.
# This rewrites the "K" column values (column G) for rows 2-56 based on the
# re-derived std/mean calculation pipeline (s_vals), replacing the old
# "Strike#" derived values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 0
    3  = 0
    4  = 1
    5  = 1
    6  = 0
    7  = 2
    8  = 0
    9  = 2
    10 = 2
    11 = 2
    12 = 1
    13 = 2
    14 = 2
    15 = 1
    16 = 2
    17 = 1
    18 = 1
    19 = 2
    20 = 2
    21 = 0
    22 = 1
    23 = 2
    24 = 2
    25 = 2
    26 = 2
    27 = 0
    28 = 0
    29 = 1
    30 = 2
    31 = 1
    32 = 0
    33 = 3
    34 = 1
    35 = 1
    36 = 2
    37 = 2
    38 = 0
    39 = 0
    40 = 2
    41 = 2
    42 = 2
    43 = 0
    44 = 0
    45 = 2
    46 = 3
    47 = 2
    48 = 3
    49 = 3
    50 = 3
    51 = 2
    52 = 1
    53 = 0
    54 = 1
    55 = 1
    56 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
